$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells that are being updated to remain Text-formatted,
# so numeric-looking strings (e.g. "0.220", "1.00", "98.565.62") keep their
# exact original text representation instead of being auto-converted to numbers.
$dCells = @("D2","D3","D5","D6","D7","D8","D9","D11","D12","D13","D15","D16","D17","D18","D19","D20","D21","D23","D24","D25","D26","D27","D28","D29","D30","D31","D34","D35","D36","D37","D38","D39","D40","D41","D43","D44","D45","D46","D48","D49","D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '98.565.62'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '3.418.76'
$ws.Range("E3").Value = '  +2.01%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '257.56'
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").Value = '668.42'
$ws.Range("E6").Value = '  +6.38%  '
$ws.Range("D7").Value = '1.53'
$ws.Range("E7").Value = '  +3.17%  '
$ws.Range("D8").Value = '0.466'
$ws.Range("E8").Value = '  +17.06%  '
$ws.Range("D9").Value = '1.06'
$ws.Range("E9").Value = '  +15.27%  '
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").Value = '3.413.53'
$ws.Range("E11").Value = '  +1.90%  '
$ws.Range("D12").Value = '0.220'
$ws.Range("E12").Value = '  +9.80%  '
$ws.Range("D13").Value = '42.71'
$ws.Range("E13").Value = '  +8.53%  '
$ws.Range("E14").Value = '  +7.93%  '
$ws.Range("D15").Value = '98.677.98'
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '4.055.25'
$ws.Range("E16").Value = '  +2.02%  '
$ws.Range("B17").Value = 'Toncoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D17").Value = '5.78'
$ws.Range("E17").Value = '  +4.02%  '
$ws.Range("D18").Value = '8.15'
$ws.Range("E18").Value = '  +31.13%  '
$ws.Range("D19").Value = '3.414.32'
$ws.Range("E19").Value = '  +1.90%  '
$ws.Range("D20").Value = '17.56'
$ws.Range("E20").Value = '  +13.12%  '
$ws.Range("D21").Value = '529.59'
$ws.Range("E21").Value = '  +8.45%  '
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("D23").Value = '10.76'
$ws.Range("E23").Value = '  +12.89%  '
$ws.Range("D24").Value = '0.0000218'
$ws.Range("E24").Value = '  +4.47%  '
$ws.Range("D25").Value = '0.435'
$ws.Range("E25").Value = '  +44.68%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = '6.34'
$ws.Range("E26").Value = '  +11.72%  '
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = '102.56'
$ws.Range("E27").Value = '  +14.54%  '
$ws.Range("D28").Value = '12.70'
$ws.Range("E28").Value = '  +5.31%  '
$ws.Range("D29").Value = '3.590.03'
$ws.Range("E29").Value = '  +1.61%  '
$ws.Range("D30").Value = '0.151'
$ws.Range("E30").Value = '  +9.57%  '
$ws.Range("D31").Value = '11.58'
$ws.Range("E31").Value = '  +17.08%  '
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = '30.60'
$ws.Range("E34").Value = '  +7.50%  '
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("D36").Value = '0.565'
$ws.Range("E36").Value = '  +19.74%  '
$ws.Range("D37").Value = '2.16'
$ws.Range("E37").Value = '  +9.92%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '0.161'
$ws.Range("E38").Value = '  +8.30%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").Value = '7.77'
$ws.Range("E39").Value = '  +5.76%  '
$ws.Range("D40").Value = '529.04'
$ws.Range("E40").Value = '  +5.45%  '
$ws.Range("D41").Value = '1.36'
$ws.Range("E41").Value = '  +6.79%  '
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("D43").Value = '3.84'
$ws.Range("E43").Value = '  +3.99%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '0.859'
$ws.Range("E44").Value = '  +4.58%  '
$ws.Range("D45").Value = '3.47'
$ws.Range("E45").Value = '  +3.61%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0425'
$ws.Range("E46").Value = '  +27.15%  '
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '8.09'
$ws.Range("E48").Value = '  +13.79%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '2.11'
$ws.Range("E49").Value = '  +7.63%  '
$ws.Range("D50").Value = '5.22'
$ws.Range("E50").Value = '  +10.61%  '
$ws.Range("E51").Value = '  +11.72%  '
